$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" conversion note text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.55 = 25871.56 pesos`n✅ 25871.56 pesos = 6.54 = 980.73 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" sheet rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 152.6
$wsTasas.Range("O10").Value = 3948
$wsTasas.Range("N12").Value = 3957
$wsTasas.Range("O12").Value = 150
